$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B7 (name of the row with id=6) from "bob 3" to "bob 32"
$ws.Range("B7").Value = "bob 32"

# Update the selected cell to match the author's last selection
$ws.Range("C10").Select() | Out-Null
